# Auto-generated Excel COM-interop script to update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.372.52"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("E3").Value = "  +4.24%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'323.59"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5101"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.4131"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "'0.08708"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").Value = "'1.136"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'24.73"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "1.996.80"
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").Value = "'6.544"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "'7.431"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "'0.9979"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'94.09"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "'0.00001115"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "'0.06508"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "'6.173"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "30.422.35"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").Value = "'2.213"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "2.229.67"
$ws.Range("E26").Value = "  +4.98%  "
$ws.Range("D27").Value = "'22.44"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "'163.24"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").Value = "'2.411"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").Value = "'131.57"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "'1.141"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'6.058"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").Value = "'3.833"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "'1.333"
$ws.Range("E35").Value = "  +12.06%  "
$ws.Range("D36").Value = "'0.02511"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "'0.06599"
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("D38").Value = "'5.373"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "'12.21"
$ws.Range("E39").Value = "  +7.19%  "
$ws.Range("D40").Value = "'0.2198"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").Value = "'9.009"
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("D42").Value = "'0.6608"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "'1.230"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'13.70"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").Value = "'0.6149"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "'2.201"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'3.664"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'1.260"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").Value = "'124.54"
$ws.Range("D50").Value = "'80.31"
$ws.Range("E51").Value = "  +1.29%  "

# Reset style on cells that needed a quote-prefix to stay numeric-looking text,
# so no extra formatting/style is left on the cell (matches original plain style).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
